# Battery Cost per Unit Cap.xlsx - "Adjust cost files relying on US data"
#
# This script:
#  1. Adds an "India:US cost adjustment" factor (with a source note) to the
#     About sheet (rows 36-38), including the numeric adjustment factor in A38.
#  2. Multiplies every formula in BCpUC!B2:B38 by the new About!$A$38 factor.
#  3. Multiplies the BBoSCpUC!B2 formula by the same About!$A$38 factor
#     (the rest of row 2 derives from B2, so it recalculates automatically).
#  4. Restores the view/selection state that results from these edits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. About sheet: add India:US cost adjustment rows
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A36").Value = "India:US cost adjustment"
$wsAbout.Range("A37").Value = "see ""scaling-factors.xlsx in the InputData folder for source information."
$wsAbout.Range("A38").Value = 0.50596615326007366

# ---------------------------------------------------------------------------
# 2. BCpUC sheet: scale every year's cost by the India:US cost adjustment
# ---------------------------------------------------------------------------
$wsBCpUC = $wb.Worksheets.Item("BCpUC")

for ($row = 2; $row -le 38; $row++) {
    $calcRow = $row + 78
    $cell = $wsBCpUC.Cells.Item($row, 2)
    $cell.Formula = "=Calculations!C$calcRow*1000*About!`$A`$38"
}

# ---------------------------------------------------------------------------
# 3. BBoSCpUC sheet: scale the balance-of-system cost by the same factor
# ---------------------------------------------------------------------------
$wsBBoSCpUC = $wb.Worksheets.Item("BBoSCpUC")
$wsBBoSCpUC.Range("B2").Formula = "='Balance of System Calculations'!D26*About!`$A`$38"

# ---------------------------------------------------------------------------
# 4. Restore view/selection state
# ---------------------------------------------------------------------------
$wsBoSCalc = $wb.Worksheets.Item("Balance of System Calculations")
$wsBoSCalc.Activate()
$wsBoSCalc.Range("F22").Select()

$wsBCpUC.Activate()
$wsBCpUC.Range("B2:B38").Select()

$wsBBoSCpUC.Activate()
$wsBBoSCpUC.Range("R2").Select()

$wsAbout.Activate()
$wsAbout.Range("A36:A38").Select()
